# Add a new "2022-Q4" sheet (with its quarterly fund-holding data) right
# after "总计" and before "2022-Q3", and record the new quarter in the
# "总计" summary sheet.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Create the new "2022-Q4" sheet by duplicating the existing
#    "2022-Q3" sheet (position 2) so it inherits the same column
#    layout / header text / cell styles (bold header, bordered index
#    column, etc.) instead of starting from a blank sheet.
# ------------------------------------------------------------------
$q3 = $wb.Worksheets.Item(2)
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# ------------------------------------------------------------------
# 2) Overwrite the copied rows with the real 2022-Q4 holdings.
#    Columns B-G are stored as text in this workbook (fund codes keep
#    leading zeros, numbers keep trailing zeros), so they are assigned
#    with a leading apostrophe to force text even when the content
#    looks numeric. Columns A (row index) and H (rank) are numeric.
# ------------------------------------------------------------------
$q4data = @(
    @("159851", "华宝中证金融科技主题ETF",  "1.98", "98.27", "4.49", "0.0889", 4),
    @("516100", "华夏中证金融科技主题ETF",  "0.60", "97.54", "4.46", "0.0268", 4),
    @("516860", "博时中证金融科技主题ETF",  "0.27", "98.47", "4.51", "0.0122", 4),
    @("002135", "广发鑫源灵活配置混合A",    "0.59", "27.74", "1.42", "0.0084", 5),
    @("002020", "国都创新驱动灵活配置混合", "0.12", "83.47", "2.98", "0.0036", 10),
    @("002136", "广发鑫源灵活配置混合C",    "0.19", "27.74", "1.42", "0.0027", 5)
)

for ($i = 0; $i -lt $q4data.Count; $i++) {
    $r = $i + 2
    $row = $q4data[$i]
    $q4.Cells.Item($r, 2).Value = "'" + $row[0]
    $q4.Cells.Item($r, 3).Value = $row[1]
    $q4.Cells.Item($r, 4).Value = "'" + $row[2]
    $q4.Cells.Item($r, 5).Value = "'" + $row[3]
    $q4.Cells.Item($r, 6).Value = "'" + $row[4]
    $q4.Cells.Item($r, 7).Value = "'" + $row[5]
    $q4.Cells.Item($r, 8).Value = $row[6]
}

# ------------------------------------------------------------------
# 3) Update the "总计" (summary) sheet: add a row for 2022-Q4
#    (holding count 6, value 0.14) and push the existing quarters
#    down by one row so the history keeps the same chronological
#    order below the new entry.
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# The table grows by one row (A1:D7 -> A1:D8); duplicate the formatting
# (bold/bordered index-column style) of the last existing row into the
# brand-new row 8 before filling in values, so the new row matches the
# look of the others instead of being left unstyled.
$summary.Range("A7:D7").Copy()
$summary.Range("A8:D8").PasteSpecial(-4122)

$summaryData = @(
    @("2022-Q4", 6, 0.14),
    @("2022-Q3", 6, 0.31),
    @("2022-Q2", 2, 0.14),
    @("2022-Q1", 4, 0.27),
    @("2021-Q4", 2, 0.1),
    @("2021-Q3", 5, 0.22),
    @("2021-Q2", 11, 0.26)
)

for ($i = 0; $i -lt $summaryData.Count; $i++) {
    $r = $i + 2
    $row = $summaryData[$i]
    $summary.Cells.Item($r, 1).Value = $i
    $summary.Cells.Item($r, 2).Value = $row[0]
    $summary.Cells.Item($r, 3).Value = $row[1]
    $summary.Cells.Item($r, 4).Value = $row[2]
}

# Keep "总计" as the active/selected sheet (as it was before the edit);
# copying Q3 into a new tab leaves the copy selected otherwise.
$summary.Activate()

